$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 8..15 (extr1..extr8) down to rows 10..17, to make room
# for two new "line7"/"line8" rows at 8..9. Process bottom-up so we never
# overwrite a source row before it has been copied. Copying whole A:E row
# ranges preserves cell styles / shared-string refs / boolean types exactly.
for ($oldRow = 15; $oldRow -ge 8; $oldRow--) {
    $newRow = $oldRow + 2
    $src = $ws.Range("A" + $oldRow + ":E" + $oldRow)
    $dst = $ws.Range("A" + $newRow + ":E" + $newRow)
    $src.Copy($dst)
    # The "index" column keeps counting up; bump it by 2 to account for the
    # two newly-inserted rows ahead of it.
    $ws.Cells.Item($newRow, 1).Value = $ws.Cells.Item($newRow, 1).Value() + 2
}

# New row 8: line7
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

# New row 9: line8
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# Give the new rows' A-column cells the same style as the rest of column A
# (bold/centered/bordered "index" style) instead of leaving them unstyled.
$ws.Range("A2").Copy($ws.Range("A8"))
$ws.Range("A2").Copy($ws.Range("A9"))
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7

# Fix the in_service flags that changed for the shifted extr rows.
$ws.Range("E10").Value = $true   # extr1: 0 -> 1
$ws.Range("E13").Value = $true   # extr4: 0 -> 1
$ws.Range("E14").Value = $false  # extr5: 1 -> 0
